$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SWKS")

# Row 4 - Inventory
$ws.Range("C4").Value = 719000000.0
$ws.Range("D4").Value = 806000000.0
$ws.Range("E4").Value = 698000000.0
$ws.Range("F4").Value = 649000000.0
$ws.Range("G4").Value = 604000000.0

# Row 14 - Accounts Payable
$ws.Range("C14").Value = 265000000.0
$ws.Range("D14").Value = 227000000.0
$ws.Range("E14").Value = 201000000.0
$ws.Range("F14").Value = 170000000.0
$ws.Range("G14").Value = 139000000.0

# Row 20 - Long Term Tax Liability (Deferred)
$ws.Range("C20").Value = -56000000.0
$ws.Range("D20").Value = -55000000.0
$ws.Range("E20").Value = -43000000.0
$ws.Range("F20").Value = -40000000.0
$ws.Range("G20").Value = -40000000.0

# Row 33 - Net Debt (B33 was blank inline string, now a number)
$ws.Range("B33").Value = -1423600000.0
